$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'Dekad start date (1st, 11th, or 21st of the month) representing the 10‑day rainfall observation period.'
$ws.Range('B3').Value = 'Administrative level code for the unit (e.g., 1 = region, 2 = district or municipality).'
$ws.Range('B4').Value = 'Numeric identifier of the administrative unit in the source dataset.'
$ws.Range('B5').Value = 'Standardized administrative unit code (P-code) for the subnational area (e.g., district/municipality).'
$ws.Range('B6').Value = 'Human-readable name of the administrative unit corresponding to the PCODE.'
$ws.Range('B7').Value = 'Name of the larger municipality or metropolitan region that contains the administrative unit.'
$ws.Range('B8').Value = 'Number of satellite grid cells used to compute rainfall metrics for the unit (indicator of spatial coverage and data quality).'
$ws.Range('B9').Value = '10‑day (dekadal) rainfall amount in millimetres for the given period.'
$ws.Range('B10').Value = 'Long-term average 10‑day rainfall (mm) for that calendar dekad at the unit (climatological mean).'
$ws.Range('B11').Value = '1‑month rolling rainfall total (mm), typically the sum of rainfall over the last three dekads.'
$ws.Range('B12').Value = 'Long-term average 1‑month rolling rainfall (mm) for the same calendar window.'
$ws.Range('B13').Value = '3‑month rolling rainfall total (mm), typically the sum of rainfall over the last nine dekads.'
$ws.Range('B14').Value = 'Long-term average 3‑month rolling rainfall (mm) for the same calendar window.'
$ws.Range('B15').Value = '10‑day rainfall anomaly in percent, representing how rfh deviates from rfh_avg.'
$ws.Range('B16').Value = '1‑month rainfall anomaly in percent, representing how r1h deviates from r1h_avg.'
$ws.Range('B17').Value = '3‑month rainfall anomaly in percent, representing how r3h deviates from r3h_avg.'
$ws.Range('B18').Value = 'Data product type label: “forecast”, “prelim” (preliminary), or “final” observation.'
$ws.Range('B19').Value = 'Previous dekad’s 10‑day rainfall (mm) for the same PCODE (lag of 1 dekad).'
$ws.Range('B20').Value = '10‑day rainfall (mm) from three dekads earlier (approximately one month before) for the same PCODE.'
$ws.Range('B21').Value = 'Previous dekad’s 1‑month rainfall anomaly (%) for the same PCODE.'
$ws.Range('B22').Value = 'Calendar month number (1–12) derived from the dekad date.'
$ws.Range('B23').Value = 'Position of the dekad within the month: 1 for days 1–10, 2 for days 11–20, 3 for days 21–end.'
$ws.Range('B24').Value = 'Indicator (0 = no, 1 = yes) that the dekad falls within the defined rainy season (e.g., May–October).'
$ws.Range('B25').Value = 'Sine transformation of the month value, used to encode cyclical seasonality for machine learning models.'
$ws.Range('B26').Value = 'Cosine transformation of the month value, used to encode cyclical seasonality for machine learning models.'
$ws.Range('B27').Value = 'Indicator (0 = no, 1 = yes) that the dekad is classified as drought-like based on a low rfq threshold.'
$ws.Range('B28').Value = 'Indicator (0 = no, 1 = yes) that the dekad is classified as having unusually high rainfall based on a high rfq threshold.'
$ws.Range('B29').Value = '6‑month moving average of the 3‑month rolling rainfall (r3h) for each PCODE, capturing medium-term trends.'
$ws.Range('B30').Value = '12‑month moving average of the 3‑month rolling rainfall (r3h) for each PCODE, capturing long-term trends.'
$ws.Range('B31').Value = 'Rolling standard deviation of 10‑day rainfall over the last three dekads, indicating short-term rainfall variability.'
$ws.Range('B32').Value = 'Change in 10‑day rainfall from the previous dekad, computed as rfh − rfh_lag_1 (mm).'

$ws.Range('B35').Select() | Out-Null

